$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.632.25"
$ws.Range("E2").Value = "  +1.84%  "

$ws.Range("D3").Value = "1.889.35"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4918"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2957"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06794"
$ws.Range("E9").Value = "  +2.68%  "

$ws.Range("D10").Value = "1.887.05"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("E11").Value = "  +3.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07240"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.28"
$ws.Range("E13").Value = "  +5.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6773"
$ws.Range("E14").Value = "  +1.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.041"
$ws.Range("E15").Value = "  +2.96%  "

$ws.Range("D16").Value = "30.610.82"
$ws.Range("E16").Value = "  +1.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007965"
$ws.Range("E17").Value = "  +1.95%  "

$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("E19").Value = "  +2.95%  "

$ws.Range("D20").Value = "2.131.21"
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("E22").Value = "  +0.96%  "

$ws.Range("E23").Value = "  +35.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.067"
$ws.Range("E24").Value = "  +3.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.336"
$ws.Range("E25").Value = "  +2.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.00"
$ws.Range("E26").Value = "  +2.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.21"
$ws.Range("E27").Value = "  +13.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.906"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.409"
$ws.Range("E29").Value = "  +1.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.331"
$ws.Range("E30").Value = "  +3.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09085"
$ws.Range("E31").Value = "  +3.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.015"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05197"
$ws.Range("E33").Value = "  +3.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7594"
$ws.Range("E34").Value = "  +4.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.110"
$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.772"
$ws.Range("E36").Value = "  +4.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01841"
$ws.Range("E37").Value = "  +1.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.673"
$ws.Range("E38").Value = "  -0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.148"
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9350"
$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4422"
$ws.Range("E41").Value = "  +4.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.35"
$ws.Range("E42").Value = "  +1.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.759"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.597"
$ws.Range("E45").Value = "  +3.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1346"
$ws.Range("E46").Value = "  +5.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05859"
$ws.Range("E47").Value = "  +2.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.710"
$ws.Range("E48").Value = "  +5.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.427"
$ws.Range("E49").Value = "  +6.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3922"
$ws.Range("E50").Value = "  +3.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.54"
$ws.Range("E51").Value = "  +2.28%  "

